$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# Row 55: date header 2023/01/09 (style copied from existing date cell A1, s="1")
$ws.Range("A1").Copy()
$ws.Range("A55").PasteSpecial($xlPasteFormats)
$ws.Range("A55").Value = 44935

# Row 56: time + text (style copied from existing time cell A3, s="2")
$ws.Range("A3").Copy()
$ws.Range("A56").PasteSpecial($xlPasteFormats)
$ws.Range("A56").Value = 0.46527777777777773
$ws.Range("B56").Value = "三層式 錯過'"

# Row 57: time + text
$ws.Range("A3").Copy()
$ws.Range("A57").PasteSpecial($xlPasteFormats)
$ws.Range("A57").Value = 0.49305555555555558
$ws.Range("B57").Value = "講故事 錯過"

# Row 58: time + text
$ws.Range("A3").Copy()
$ws.Range("A58").PasteSpecial($xlPasteFormats)
$ws.Range("A58").Value = 0.6791666666666667
$ws.Range("B58").Value = "EditUser"

$excel.CutCopyMode = 0

$ws.Range("B59").Select()
